{"js": "// Renumber the \"_Toc...\" bookmarks left behind by a stale Word TOC field\n// update. Each bookmark keeps its position/id \u2014 only the w:name changes.\nconst renames = [\n  [\"_Toc5728479\", \"_Toc5791558\"],\n  [\"_Toc5728480\", \"_Toc5791559\"],\n  [\"_Toc5728481\", \"_Toc5791560\"],\n  [\"_Toc5728482\", \"_Toc5791561\"],\n  [\"_Toc5728483\", \"_Toc5791562\"],\n  [\"_Toc5728484\", \"_Toc5791563\"],\n];\n\nfor (const [oldName, newName] of renames) {\n  const range = context.document.getBookmarkRange(oldName);\n  context.document.deleteBookmark(oldName);\n  range.insertBookmark(newName);\n}\n\nawait context.sync();\n", "ps1": "# Renumber the \"_Toc...\" bookmarks left behind by a stale Word TOC field\n# update. Each bookmark keeps its position/id - only the name changes, so\n# we add a new bookmark with the new name over the old bookmark's range,\n# then delete the old one (Bookmark.Name is read-only in the Word object\n# model, so a direct rename is not available).\n$d = $word.ActiveDocument\n\n$renames = @(\n    @{ Old = \"_Toc5728479\"; New = \"_Toc5791558\" },\n    @{ Old = \"_Toc5728480\"; New = \"_Toc5791559\" },\n    @{ Old = \"_Toc5728481\"; New = \"_Toc5791560\" },\n    @{ Old = \"_Toc5728482\"; New = \"_Toc5791561\" },\n    @{ Old = \"_Toc5728483\"; New = \"_Toc5791562\" },\n    @{ Old = \"_Toc5728484\"; New = \"_Toc5791563\" }\n)\n\nforeach ($pair in $renames) {\n    $bm = $d.Bookmarks($pair.Old)\n    $r = $bm.Range\n    $d.Bookmarks.Add($pair.New, $r)\n    $bm.Delete()\n}\n"}
